$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Pedro"
$ws.Cells.Item(4, 3).Value = "Jos"
$ws.Cells.Item(4, 4).Value = "Armas"
$ws.Cells.Item(4, 5).Value = "Coyo"

$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "1234455"
$ws.Cells.Item(4, 6).Style = "Normal"

$ws.Cells.Item(4, 7).Value = "jospspd"
$ws.Cells.Item(4, 8).Value = "kaddkdakakd"

$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = "1987-12-01"
$ws.Cells.Item(4, 9).Style = "Normal"
